# Auto-generated script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '34.596.32'
$ws.Cells.Item(3, 4).Value = '1.803.19'
$ws.Cells.Item(3, 5).Value = '  +0.94%  '
$ws.Cells.Item(4, 5).Value = '  -0.12%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '224.21'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.45%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.552'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.22%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.10%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '32.48'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +3.53%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.289'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +2.89%  '
$ws.Cells.Item(10, 5).Value = '  +7.96%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0929'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.02%  '
$ws.Cells.Item(12, 4).Value = '2.063.38'
$ws.Cells.Item(12, 5).Value = '  +0.95%  '
$ws.Cells.Item(13, 2).Value = 'Chainlink'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.09'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -3.66%  '
$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value = '1.798.85'
$ws.Cells.Item(14, 5).Value = '  +0.79%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.642'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.30%  '
$ws.Cells.Item(16, 4).Value = '34.626.25'
$ws.Cells.Item(16, 5).Value = '  +1.63%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.32'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +2.07%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '69.19'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.48%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '252.46'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.17%  '
$ws.Cells.Item(20, 5).Value = '  +8.19%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.05'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +5.72%  '
$ws.Cells.Item(22, 5).Value = '  -0.16%  '
$ws.Cells.Item(23, 5).Value = '  -0.31%  '
$ws.Cells.Item(24, 5).Value = '  +1.30%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '161.64'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +2.83%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.41'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -1.01%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.14'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +1.65%  '
$ws.Cells.Item(28, 5).Value = '  -0.10%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.14%  '
$ws.Cells.Item(30, 2).Value = 'Swop.fi'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/yrCr2HW2c+swopfi-swop'
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '571.53'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +994.98%  '
$ws.Cells.Item(31, 2).Value = 'Hedera'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0527'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +2.06%  '
$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.80'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.20%  '
$ws.Cells.Item(33, 2).Value = 'PancakeSwap'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.20'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.21%  '
$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.62'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.40%  '
$ws.Cells.Item(35, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.88'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +2.67%  '
$ws.Cells.Item(36, 2).Value = 'Maker'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(36, 4).Value = '1.433.34'
$ws.Cells.Item(36, 5).Value = '  -1.17%  '
$ws.Cells.Item(37, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.07'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.08%  '
$ws.Cells.Item(38, 2).Value = 'ImmutableX'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.643'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +2.38%  '
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0192'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +2.67%  '
$ws.Cells.Item(40, 2).Value = 'Aave'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '84.72'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +1.44%  '
$ws.Cells.Item(41, 2).Value = 'ARBITRUM'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.961'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +6.78%  '
$ws.Cells.Item(42, 2).Value = 'MXToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.80'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.69%  '
$ws.Cells.Item(43, 2).Value = 'HuobiToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.35'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.02%  '
$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.16'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +3.95%  '
$ws.Cells.Item(45, 2).Value = 'FraxShare'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.04'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +4.49%  '
$ws.Cells.Item(46, 2).Value = 'WEMIXToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.06'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.97%  '
$ws.Cells.Item(47, 2).Value = 'Kaspa'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0498'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -2.67%  '
$ws.Cells.Item(48, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(48, 4).Value = '1.957.16'
$ws.Cells.Item(48, 5).Value = '  +0.65%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '106.78'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +9.14%  '
$ws.Cells.Item(50, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.31'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +3.61%  '
$ws.Cells.Item(51, 2).Value = 'PaxDollar'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.02%  '
